$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 25.99000000000062
$ws.Range("H2").Value = [double]"6.254777603522009e-16"
$ws.Range("K2").Value = 51.27483781071916
$ws.Range("L2").Value = "[42.18173779719306, 60.36793782424527]"
$ws.Range("O2").Value = 1.792500312859041
$ws.Range("P2").Value = "[1.6038160694001942, 1.9811845563178885]"
$ws.Range("S2").Value = 61.11993207957558
$ws.Range("T2").Value = "[55.73809853290837, 66.5017656262428]"
$ws.Range("W2").Value = 18.57543543543588
$ws.Range("X2").Value = 17.79495495495538
$ws.Range("Y2").Value = 19.35591591591638

$ws.Range("E3").Value = 25.68000000000058
$ws.Range("G3").Value = [double]"3.129185799366496e-11"
$ws.Range("H3").Value = [double]"1.147234159522744e-10"
$ws.Range("K3").Value = 46.18403522373059
$ws.Range("L3").Value = "[30.20484298903817, 62.16322745842301]"
$ws.Range("M3").Value = [double]"6.315741507556538e-08"
$ws.Range("N3").Value = [double]"6.315741507556538e-08"
$ws.Range("O3").Value = 2.106974051957119
$ws.Range("P3").Value = "[1.7421845146033492, 2.4717635893108882]"
$ws.Range("S3").Value = 65.00901643433575
$ws.Range("T3").Value = "[56.25625669025756, 73.76177617841394]"
$ws.Range("W3").Value = 17.06858858858897
$ws.Range("X3").Value = 15.57765765765801
$ws.Range("Y3").Value = 18.55951951951993
